$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the col_value column (C) for rows 3-9: wrap the existing numeric-looking
# text in literal single quotes (e.g. 10 -> '0010', 0 -> '0000', 253 -> '253', 5511 -> '5511').
# Typing a value that begins with two apostrophes makes Excel store the text with
# a single leading literal apostrophe (quote-prefix applied) while keeping the
# trailing apostrophe as literal text too, matching stored cell text like '0010'.
$ws.Range("C3").Value = "''0010'"
$ws.Range("C4").Value = "''0000'"
$ws.Range("C5").Value = "''0000'"
$ws.Range("C6").Value = "''253'"
$ws.Range("C7").Value = "''5511'"
$ws.Range("C8").Value = "''0010'"
$ws.Range("C9").Value = "''0000'"

# Row 2 (the long wrapped Formula text) is re-measured by Excel after the edits
# above and its autofit height settles lower than before (224 -> 208).
$ws.Rows.Item(2).RowHeight = 208

# Move / update the active selection to reflect the last edited cell.
$ws.Range("C8").Select()
